# "fix(gui) step 1 and 2"
#
# The price list's heading date is advanced by one day, and the unit
# prices for the two "MIRADOR OPTICO" items (MO-100 / step 1, and
# MO-101 / step 2) are updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the heading date in A1 (45308 -> 45309, i.e. one day later) ---
$ws.Range("A1").Value = 45309

# --- Step 1: MO-100 "MIRADOR OPTICO Hº Niquel." price ---
$ws.Range("D22").Value = 2950.798

# --- Step 2: MO-101 "MIRADOR OPTICO Hº Bceado." price ---
$ws.Range("D23").Value = 2950.798
